$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 103/104/105: cyclic rotation of match details (F:V), A:E (index/date) unchanged ---
$v103 = $ws.Range("F103:V103").Value()
$v104 = $ws.Range("F104:V104").Value()
$v105 = $ws.Range("F105:V105").Value()
$ws.Range("F103:V103").Value = $v104
$ws.Range("F104:V104").Value = $v105
$ws.Range("F105:V105").Value = $v103

# --- Row 230/231: swap match details (F:V) ---
$v230 = $ws.Range("F230:V230").Value()
$v231 = $ws.Range("F231:V231").Value()
$ws.Range("F230:V230").Value = $v231
$ws.Range("F231:V231").Value = $v230

# --- Row 237/238: swap match details (F:V) ---
$v237 = $ws.Range("F237:V237").Value()
$v238 = $ws.Range("F238:V238").Value()
$ws.Range("F237:V237").Value = $v238
$ws.Range("F238:V238").Value = $v237

# --- New rows 245-250, appended after existing row 244, formatted like row 244 ---

$ws.Range("A244:V244").Copy()
$ws.Range("A245:V245").PasteSpecial(-4122)
$ws.Range("A245").Value = 244
$ws.Range("B245").Value = 'japan'
$ws.Range("C245").Value = 'j1-league'
$ws.Range("D245").Value = "'2023"
$ws.Range("E245").Value = 45192.29166666666
$ws.Range("F245").Value = 'Albirex Niigata'
$ws.Range("G245").Value = 3
$ws.Range("H245").Value = 'Yokohama FC'
$ws.Range("I245").Value = 1
$ws.Range("J245").Value = 1.72
$ws.Range("K245").Value = '17/09/2023 11:12'
$ws.Range("L245").Value = 1.92
$ws.Range("M245").Value = '23/09/2023 06:35'
$ws.Range("N245").Value = 3.98
$ws.Range("O245").Value = '17/09/2023 11:12'
$ws.Range("P245").Value = 3.64
$ws.Range("Q245").Value = '23/09/2023 06:35'
$ws.Range("R245").Value = 4.98
$ws.Range("S245").Value = '17/09/2023 11:12'
$ws.Range("T245").Value = 4.31
$ws.Range("U245").Value = '23/09/2023 06:35'
$ws.Range("V245").Value = 'https://www.betexplorer.com/football/japan/j1-league/albirex-niigata-yokohama-fc/QDjnLVvn/'

$ws.Range("A245:V245").Copy()
$ws.Range("A246:V246").PasteSpecial(-4122)
$ws.Range("A246").Value = 245
$ws.Range("B246").Value = 'japan'
$ws.Range("C246").Value = 'j1-league'
$ws.Range("D246").Value = "'2023"
$ws.Range("E246").Value = 45192.33333333334
$ws.Range("F246").Value = 'FC Tokyo'
$ws.Range("G246").Value = 3
$ws.Range("H246").Value = 'Sagan Tosu'
$ws.Range("I246").Value = 2
$ws.Range("J246").Value = 1.88
$ws.Range("K246").Value = '15/09/2023 11:12'
$ws.Range("L246").Value = 1.97
$ws.Range("M246").Value = '23/09/2023 07:55'
$ws.Range("N246").Value = 3.79
$ws.Range("O246").Value = '15/09/2023 11:12'
$ws.Range("P246").Value = 3.92
$ws.Range("Q246").Value = '23/09/2023 07:55'
$ws.Range("R246").Value = 4.19
$ws.Range("S246").Value = '15/09/2023 11:12'
$ws.Range("T246").Value = 3.76
$ws.Range("U246").Value = '23/09/2023 07:59'
$ws.Range("V246").Value = 'https://www.betexplorer.com/football/japan/j1-league/fc-tokyo-sagan-tosu/EZMT2lPH/'

$ws.Range("A246:V246").Copy()
$ws.Range("A247:V247").PasteSpecial(-4122)
$ws.Range("A247").Value = 246
$ws.Range("B247").Value = 'japan'
$ws.Range("C247").Value = 'j1-league'
$ws.Range("D247").Value = "'2023"
$ws.Range("E247").Value = 45192.375
$ws.Range("F247").Value = 'Nagoya Grampus'
$ws.Range("G247").Value = 1
$ws.Range("H247").Value = 'Hokkaido Consadole Sapporo'
$ws.Range("I247").Value = 1
$ws.Range("J247").Value = 2.15
$ws.Range("K247").Value = '16/09/2023 11:13'
$ws.Range("L247").Value = 2.26
$ws.Range("M247").Value = '23/09/2023 08:59'
$ws.Range("N247").Value = 3.68
$ws.Range("O247").Value = '16/09/2023 11:13'
$ws.Range("P247").Value = 3.72
$ws.Range("Q247").Value = '23/09/2023 08:59'
$ws.Range("R247").Value = 3.31
$ws.Range("S247").Value = '16/09/2023 11:13'
$ws.Range("T247").Value = 3.16
$ws.Range("U247").Value = '23/09/2023 08:59'
$ws.Range("V247").Value = 'https://www.betexplorer.com/football/japan/j1-league/nagoya-grampus-hokkaido-consadole-sapporo/Kl5kKkfh/'

$ws.Range("A247:V247").Copy()
$ws.Range("A248:V248").PasteSpecial(-4122)
$ws.Range("A248").Value = 247
$ws.Range("B248").Value = 'japan'
$ws.Range("C248").Value = 'j1-league'
$ws.Range("D248").Value = "'2023"
$ws.Range("E248").Value = 45192.5
$ws.Range("F248").Value = 'Kashiwa Reysol'
$ws.Range("G248").Value = 1
$ws.Range("H248").Value = 'Avispa Fukuoka'
$ws.Range("I248").Value = 3
$ws.Range("J248").Value = 2.28
$ws.Range("K248").Value = '17/09/2023 10:42'
$ws.Range("L248").Value = 2.48
$ws.Range("M248").Value = '23/09/2023 11:57'
$ws.Range("N248").Value = 3.19
$ws.Range("O248").Value = '17/09/2023 10:42'
$ws.Range("P248").Value = 2.82
$ws.Range("Q248").Value = '23/09/2023 11:58'
$ws.Range("R248").Value = 3.58
$ws.Range("S248").Value = '17/09/2023 10:42'
$ws.Range("T248").Value = 3.7
$ws.Range("U248").Value = '23/09/2023 11:57'
$ws.Range("V248").Value = 'https://www.betexplorer.com/football/japan/j1-league/kashiwa-reysol-avispa-fukuoka/MaQy1SgU/'

$ws.Range("A248:V248").Copy()
$ws.Range("A249:V249").PasteSpecial(-4122)
$ws.Range("A249").Value = 248
$ws.Range("B249").Value = 'japan'
$ws.Range("C249").Value = 'j1-league'
$ws.Range("D249").Value = "'2023"
$ws.Range("E249").Value = 45192.5
$ws.Range("F249").Value = 'Kyoto'
$ws.Range("G249").Value = 1
$ws.Range("H249").Value = 'Sanfrecce Hiroshima'
$ws.Range("I249").Value = 0
$ws.Range("J249").Value = 3.73
$ws.Range("K249").Value = '16/09/2023 11:13'
$ws.Range("L249").Value = 4.92
$ws.Range("M249").Value = '23/09/2023 11:55'
$ws.Range("N249").Value = 3.59
$ws.Range("O249").Value = '16/09/2023 11:13'
$ws.Range("P249").Value = 3.95
$ws.Range("Q249").Value = '23/09/2023 11:56'
$ws.Range("R249").Value = 2.06
$ws.Range("S249").Value = '16/09/2023 11:13'
$ws.Range("T249").Value = 1.74
$ws.Range("U249").Value = '23/09/2023 11:56'
$ws.Range("V249").Value = 'https://www.betexplorer.com/football/japan/j1-league/kyoto-sanfrecce-hiroshima/bu4gJ99b/'

$ws.Range("A249:V249").Copy()
$ws.Range("A250:V250").PasteSpecial(-4122)
$ws.Range("A250").Value = 249
$ws.Range("B250").Value = 'japan'
$ws.Range("C250").Value = 'j1-league'
$ws.Range("D250").Value = "'2023"
$ws.Range("E250").Value = 45192.5
$ws.Range("F250").Value = 'Vissel Kobe'
$ws.Range("G250").Value = 1
$ws.Range("H250").Value = 'Cerezo Osaka'
$ws.Range("I250").Value = 0
$ws.Range("J250").Value = 1.94
$ws.Range("K250").Value = '16/09/2023 11:13'
$ws.Range("L250").Value = 2.42
$ws.Range("M250").Value = '23/09/2023 11:57'
$ws.Range("N250").Value = 3.74
$ws.Range("O250").Value = '16/09/2023 11:13'
$ws.Range("P250").Value = 3.48
$ws.Range("Q250").Value = '23/09/2023 11:54'
$ws.Range("R250").Value = 3.96
$ws.Range("S250").Value = '16/09/2023 11:13'
$ws.Range("T250").Value = 3.06
$ws.Range("U250").Value = '23/09/2023 11:57'
$ws.Range("V250").Value = 'https://www.betexplorer.com/football/japan/j1-league/vissel-kobe-cerezo-osaka/f7c2HmvB/'

Write-Host "done"
